# Remove the whole bullet-list paragraph "Formater ou valider: téléphone "
# (the requirement about formatting/validating the phone number is dropped
# entirely per the commit "Enlegistlement: Waridation de téléphon").
$d = $word.ActiveDocument

$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "Formater\s*ou valider") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $target.Range.Delete()
}
